# Daily attendance processing - 2026-01-03 04:23:05
#
# Normalizes the "Recorded By" column (G) on the active sheet: whenever the
# first name in the comma-separated list of recorders is some case variant
# of "System", that entry is rotated to the end of the list instead of the
# front (e.g. "System, dnasr281@gmail.com" -> "dnasr281@gmail.com, System").
# Entries that don't start with "System" are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$changed = 0

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # column G = "Recorded By"
    $text = $cell.Text

    if ([string]::IsNullOrEmpty($text)) {
        continue
    }

    $parts = $text -split ", "

    if ($parts.Length -gt 1 -and $parts[0].Trim().ToLower() -eq "system") {
        $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
        $newText = $rotated -join ", "
        $cell.Value = $newText
        $changed++
    }
}

Write-Output ("Rotated 'System' to end of Recorded-By list on " + $changed + " rows")
